$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-31) holds the "Förändrad" date as an Excel serial date.
# Update the serial value from 45170 to 45174 (2023-09-01 -> 2023-09-05)
# while preserving the existing cell style/format.
for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45170) {
        $cell.Value2 = 45174
    }
}
